$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 112, shifting existing rows 112:227 down to 113:228.
$ws.Rows("112:112").Insert()

# Populate the newly inserted row 112 with the new data record.
$ws.Range("A112").Value = 9
$ws.Range("B112").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C112").Value = "Metropolitana"
$ws.Range("D112").Value = 44512
$ws.Range("D112").NumberFormat = $ws.Range("D113").NumberFormat
$ws.Range("E112").Value = 13
$ws.Range("F112").Value = 100112052
$ws.Range("G112").Value = "Albahaca"
$ws.Range("H112").Value = "Sin especificar"
$ws.Range("I112").Value = "Primera"
$ws.Range("J112").Value = 79
$ws.Range("K112").Value = 5000
$ws.Range("L112").Value = 6000
$ws.Range("M112").Value = 5494
$ws.Range("N112").Value = "$/docena de matas"
$ws.Range("O112").Value = "Provincia de Chacabuco"
$ws.Range("P112").Value = 916
$ws.Range("Q112").Value = 6
$ws.Range("R112").Value = "Hortaliza"
